$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 76 (ALC) - hunk 0
$ws.Range("H76").Value = 36143.266
$ws.Range("I76").Value = 41261.46
$ws.Range("J76").Value = 2875
$ws.Range("K76").Value = 41261.46
$ws.Range("L76").Value = 2875
$ws.Range("M76").Value = -40946.46
$ws.Range("N76").Value = -3505

# Row 79 (ALC) - hunk 1
$ws.Range("H79").Value = 36143.266
$ws.Range("I79").Value = 41261.46
$ws.Range("J79").Value = 2875
$ws.Range("K79").Value = 41261.46
$ws.Range("L79").Value = 2875
$ws.Range("M79").Value = -40169.46
$ws.Range("N79").Value = -5059

# Row 107 (ALC) - hunk 2
$ws.Range("H107").Value = 601
$ws.Range("I107").Value = 541.8570999999999
$ws.Range("J107").Value = 808
$ws.Range("K107").Value = 541.8570999999999
$ws.Range("L107").Value = 808
$ws.Range("M107").Value = 1378.1429
$ws.Range("N107").ClearContents()

# Row 113 (ALC) - hunk 3
$ws.Range("H113").Value = 1741.2084
$ws.Range("I113").Value = 1513.762
$ws.Range("J113").Value = 3333.3333
$ws.Range("K113").Value = 1513.762
$ws.Range("L113").Value = 3333.3333
$ws.Range("M113").Value = 1740.238
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 74 (ARM) - hunk 4
$ws.Range("H74").Value = 1704.2941
$ws.Range("I74").Value = 1640.9286
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 1640.9286
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -766.9286
$ws.Range("N74").Value = -3748

# Row 77 (ARM) - hunk 5
$ws.Range("H77").Value = 1704.2941
$ws.Range("I77").Value = 1640.9286
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 8204.643
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -3836.643
$ws.Range("N77").Value = -18736

# Row 80 (ARM) - hunk 6
$ws.Range("H80").Value = 27766.666
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 27766.666
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 27766.666
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -29762.666

# Row 83 (ARM) - hunk 7
$ws.Range("H83").Value = 27766.666
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 27766.666
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 83299.99800000001
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -93283.99800000001

# Row 102 (ARM) - hunk 8
$ws.Range("H102").Value = 1603.3334
$ws.Range("I102").Value = 1603.3334
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1603.3334
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 18.66660000000002

$ws = $wb.Worksheets.Item("BSM")
# Row 41 (BSM) - hunk 9
$ws.Range("H41").Value = 150163
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 150163
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 150163
$ws.Range("N41").Value = -150939

# Row 42 (BSM) - hunk 10
$ws.Range("H42").Value = 175340
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 175340
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 175340
$ws.Range("N42").Value = -175996

# Row 47 (BSM) - hunk 11
$ws.Range("H47").Value = 125342
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 125342
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 125342
$ws.Range("N47").Value = -126382

# Row 48 (BSM) - hunk 12
$ws.Range("H48").Value = 150166.33
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 150166.33
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 150166.33
$ws.Range("N48").Value = -150996.33

# Row 107 (BSM) - hunk 13
$ws.Range("H107").Value = 1544.4286
$ws.Range("I107").Value = 962.2
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 962.2
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = 957.8
$ws.Range("N107").Value = -6840

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (CRP) - hunk 14
$ws.Range("H16").Value = 1350.9
$ws.Range("I16").Value = 945.44446
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 945.44446
$ws.Range("L16").Value = 5000
$ws.Range("M16").Value = -658.44446
$ws.Range("N16").Value = -5574

# Row 62 (CRP) - hunk 15
$ws.Range("H62").Value = 4999.2
$ws.Range("I62").Value = 6250
$ws.Range("J62").Value = 4165.3335
$ws.Range("K62").Value = 6250
$ws.Range("L62").Value = 4165.3335
$ws.Range("M62").Value = -5626
$ws.Range("N62").Value = -5413.3335

# Row 65 (CRP) - hunk 16
$ws.Range("H65").Value = 4999.2
$ws.Range("I65").Value = 6250
$ws.Range("J65").Value = 4165.3335
$ws.Range("K65").Value = 31250
$ws.Range("L65").Value = 20826.6675
$ws.Range("M65").Value = -28130
$ws.Range("N65").Value = -27066.6675

# Row 113 (CRP) - hunk 17
$ws.Range("H113").Value = 1350.9
$ws.Range("I113").Value = 945.44446
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 945.44446
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = 1224.55554
$ws.Range("N113").Value = -9340

$ws = $wb.Worksheets.Item("CUL")
# Row 136 (CUL) - hunk 18
$ws.Range("H136").Value = 2037.825
$ws.Range("I136").Value = 2110
$ws.Range("J136").Value = 2027.5143
$ws.Range("K136").Value = 6330
$ws.Range("L136").Value = 6082.5429
$ws.Range("M136").Value = -1230
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (GSM) - hunk 19
$ws.Range("H80").Value = 8076.154
$ws.Range("I80").Value = 3998.3333
$ws.Range("J80").Value = 11571.429
$ws.Range("K80").Value = 3998.3333
$ws.Range("L80").Value = 11571.429
$ws.Range("M80").Value = -3000.3333
$ws.Range("N80").Value = -13567.429

# Row 83 (GSM) - hunk 20
$ws.Range("H83").Value = 8076.154
$ws.Range("I83").Value = 3998.3333
$ws.Range("J83").Value = 11571.429
$ws.Range("K83").Value = 19991.6665
$ws.Range("L83").Value = 57857.145
$ws.Range("M83").Value = -14999.6665
$ws.Range("N83").Value = -67841.145

# Row 107 (GSM) - hunk 21
$ws.Range("H107").Value = 707.8570999999999
$ws.Range("I107").Value = 630.3125
$ws.Range("J107").Value = 956
$ws.Range("K107").Value = 630.3125
$ws.Range("L107").Value = 956
$ws.Range("M107").Value = 1289.6875
$ws.Range("N107").Value = -4796

# Row 122 (GSM) - hunk 22
$ws.Range("H122").Value = 2150.1365
$ws.Range("I122").Value = 1544.1111
$ws.Range("J122").Value = 2569.6924
$ws.Range("K122").Value = 4632.3333
$ws.Range("L122").Value = 7709.0772
$ws.Range("M122").Value = -2182.3333
$ws.Range("N122").Value = -12609.0772

# Row 126 (GSM) - hunk 23
$ws.Range("H126").Value = 2783612
$ws.Range("I126").Value = 8503
$ws.Range("J126").Value = 4171166.5
$ws.Range("K126").Value = 25509
$ws.Range("L126").Value = 12513499.5
$ws.Range("M126").Value = -23039
$ws.Range("N126").Value = -12518439.5

$ws = $wb.Worksheets.Item("LTW")
# Row 81 (LTW) - hunk 24
$ws.Range("H81").Value = 39750
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 39750
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 39750
$ws.Range("N81").Value = -41746

# Row 84 (LTW) - hunk 25
$ws.Range("H84").Value = 39750
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 39750
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 119250
$ws.Range("N84").Value = -129234

$ws = $wb.Worksheets.Item("WVR")
# Row 75 (WVR) - hunk 26
$ws.Range("H75").Value = 31333.334
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 31333.334
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 31333.334
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -33205.334

# Row 78 (WVR) - hunk 27
$ws.Range("H78").Value = 31333.334
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 31333.334
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 94000.00199999999
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -103360.002
